$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Stage copies of each of the 4 titled column-blocks (title row 4 through scratch rows 18-30)
$ws.Range("G4:I30").Copy($ws.Range("AZ4"))
$ws.Range("K4:M30").Copy($ws.Range("BC4"))
$ws.Range("O4:Q30").Copy($ws.Range("BF4"))
$ws.Range("S4:U30").Copy($ws.Range("BI4"))

# Paste rotated one block to the left: new G = old S, new K = old G, new O = old K, new S = old O
$ws.Range("G4:I30").Clear()
$ws.Range("BI4:BK30").Copy($ws.Range("G4"))

$ws.Range("K4:M30").Clear()
$ws.Range("AZ4:BB30").Copy($ws.Range("K4"))

$ws.Range("O4:Q30").Clear()
$ws.Range("BC4:BE30").Copy($ws.Range("O4"))

$ws.Range("S4:U30").Clear()
$ws.Range("BF4:BH30").Copy($ws.Range("S4"))

# Clear staging
$ws.Range("AZ4:BK30").Clear()
